# Adds a new data row (row 11) to the active worksheet, replicating a
# species-observation record (Dryocopus martius / Spillkraka) that was
# appended after the existing last row (row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 11

# --- Numeric columns -------------------------------------------------
$ws.Cells.Item($row, 1).Value  = 131289566      # A  Id
$ws.Cells.Item($row, 2).Value  = 57881           # B  Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value  = 100049          # E  TaxonId
$ws.Cells.Item($row, 17).Value = 567450          # Q  Ost
$ws.Cells.Item($row, 18).Value = 6510180         # R  Nord
$ws.Cells.Item($row, 19).Value = 10              # S  Noggrannhet

# --- Text columns ------------------------------------------------------
$ws.Cells.Item($row, 4).Value  = "NT"                              # D Rödlistade
$ws.Cells.Item($row, 6).Value  = "Spillkråka"                       # F Artnamn
$ws.Cells.Item($row, 7).Value  = "Dryocopus martius"                # G Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(Linnaeus, 1758)"                 # H Auktor
$ws.Cells.Item($row, 13).Value = "gammalt bo"                       # M Aktivitet
$ws.Cells.Item($row, 16).Value = "Strax öster om Björkliden, Ög"    # P Lokalnamn
$ws.Cells.Item($row, 20).Value = "Östergötland"                     # T Län
$ws.Cells.Item($row, 21).Value = "Norrköping"                       # U Kommun
$ws.Cells.Item($row, 22).Value = "Östergötland"                     # V Provins
$ws.Cells.Item($row, 23).Value = "Simonstorp"                       # W Socken
$ws.Cells.Item($row, 49).Value = "Anette Källman"                   # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Anette Källman"                   # AX Observatörer

# --- Date-looking values that must stay as plain text ------------------
# (use a leading apostrophe so Excel keeps them as text, not a date serial)
$ws.Cells.Item($row, 25).Value = "'2026-02-21"   # Y  Startdatum
$ws.Cells.Item($row, 27).Value = "'2026-02-21"   # AA Slutdatum

# --- Boolean columns -----------------------------------------------------
$ws.Cells.Item($row, 30).Value = $false          # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false          # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false          # AG Ospontan

Write-Host "Row 11 populated"
